$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph (the 2nd paragraph,
#    right after the H1 title paragraph).
$metaPara = $d.Paragraphs.Item(2)
[void]$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Desperate Dawgs 2 Gigablox for Free |
#    Review" right before the final paragraph (which currently holds the
#    "Create a cartoon-style feature image..." image-prompt text).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
[void]$lastPara.Range.InsertParagraphBefore()

$count2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count2 - 1)
[void]$newPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Desperate Dawgs 2 Gigablox for Free | Review</w:t></w:r></w:p>")

# 3. Replace the text of the final paragraph (still italic) with the new
#    meta-description copy, leaving the italic run formatting intact.
$count3 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count3)
$fr = $finalPara.Range
$frTrimmed = $d.Range($fr.Start, $fr.End - 1)
$frTrimmed.Text = "Read our review of Desperate Dawgs 2 Gigablox, a Wild West-themed slot game with Gigablox mechanic and 3 bonus features. Play for free here!"
